# Atualizado por script em 23-11-2023 14:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 150 and 151 had their match data (columns F..V) swapped.
#    Columns A..E (index/pais/torneio/temporada/data_partida) stay put.
# ---------------------------------------------------------------------
$row150 = @($ws.Range("F150").Value2, $ws.Range("G150").Value2, $ws.Range("H150").Value2, $ws.Range("I150").Value2, `
            $ws.Range("J150").Value2, $ws.Range("K150").Value2, $ws.Range("L150").Value2, $ws.Range("M150").Value2, `
            $ws.Range("N150").Value2, $ws.Range("O150").Value2, $ws.Range("P150").Value2, $ws.Range("Q150").Value2, `
            $ws.Range("R150").Value2, $ws.Range("S150").Value2, $ws.Range("T150").Value2, $ws.Range("U150").Value2, `
            $ws.Range("V150").Value2)

$row151 = @($ws.Range("F151").Value2, $ws.Range("G151").Value2, $ws.Range("H151").Value2, $ws.Range("I151").Value2, `
            $ws.Range("J151").Value2, $ws.Range("K151").Value2, $ws.Range("L151").Value2, $ws.Range("M151").Value2, `
            $ws.Range("N151").Value2, $ws.Range("O151").Value2, $ws.Range("P151").Value2, $ws.Range("Q151").Value2, `
            $ws.Range("R151").Value2, $ws.Range("S151").Value2, $ws.Range("T151").Value2, $ws.Range("U151").Value2, `
            $ws.Range("V151").Value2)

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "150").Value = $row151[$i]
    $ws.Range($cols[$i] + "151").Value = $row150[$i]
}

# ---------------------------------------------------------------------
# 2) Append a new match as row 155, copying the formatting from row 154
#    (bold/bordered index cell style + datetime style for data_partida).
# ---------------------------------------------------------------------
$ws.Range("A154:E154").Copy() | Out-Null
$ws.Range("A155:E155").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A155").Value = 154
$ws.Range("B155").Value = "colombia"
$ws.Range("C155").Value = "primera-b"

# D155 must stay a text string ("2023"), not get auto-coerced to a number.
# Temporarily force text format, enter the value, then restore the original
# (unstyled/General) formatting by re-pasting it from a plain text neighbour.
$ws.Range("D155").NumberFormat = "@"
$ws.Range("D155").Value = "2023"
$ws.Range("C154").Copy() | Out-Null
$ws.Range("D155").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E155").Value = 45253.0625
$ws.Range("F155").Value = "Patriotas"
$ws.Range("G155").Value = 3
$ws.Range("H155").Value = "Fortaleza"
$ws.Range("I155").Value = 1
$ws.Range("J155").Value = 2.56
$ws.Range("K155").Value = "20/11/2023 09:12"
$ws.Range("L155").Value = 2.61
$ws.Range("M155").Value = "23/11/2023 01:29"
$ws.Range("N155").Value = 3
$ws.Range("O155").Value = "20/11/2023 09:12"
$ws.Range("P155").Value = 2.73
$ws.Range("Q155").Value = "23/11/2023 01:07"
$ws.Range("R155").Value = 3.03
$ws.Range("S155").Value = "20/11/2023 09:12"
$ws.Range("T155").Value = 3.4
$ws.Range("U155").Value = "23/11/2023 01:29"
$ws.Range("V155").Value = "https://www.betexplorer.com/football/colombia/primera-b/patriotas-fortaleza-c-e-i-f/QXjjCEaU/"

Write-Host "done"
